# Updated the workflow for uploading renamed files to SharePoint
# Replace the old "FL Renaming / FLOBOT" paths with the new
# "WAHP Matching Automation / WABOT" paths in both config sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Config" --------------------------------------------------
$wsConfig = $wb.Worksheets.Item("Config")

# DeleteDirectory
$wsConfig.Range("B2").Value = "E:\Bot_Files\RPA WAHP Matching Automation\WABOT"
# MasterFolder
$wsConfig.Range("B6").Value = "E:\Bot_Files\RPA WAHP Matching Automation\WABOT\"
# ZippedDirectory
$wsConfig.Range("B8").Value = "E:\Bot_Files\RPA WAHP Matching Automation\WABOT"

# --- Sheet "ConfigOptions" --------------------------------------------
$wsOptions = $wb.Worksheets.Item("ConfigOptions")

# Section "WAHP SP to Local Computer - Testing"
# DeleteDirectory
$wsOptions.Range("B17").Value = "C:\Users\RollLe01\OneDrive - Reed Elsevier Group ICO Reed Elsevier Inc\Desktop\WAHP"
# MasterFolder
$wsOptions.Range("B21").Value = "C:\Users\RollLe01\OneDrive - Reed Elsevier Group ICO Reed Elsevier Inc\Desktop\WAHP\"
# ZippedDirectory
$wsOptions.Range("B23").Value = "C:\Users\RollLe01\OneDrive - Reed Elsevier Group ICO Reed Elsevier Inc\Desktop\WAHP"

# Section "WAHP SP to Virtual Desktop - Testing"
# MasterFolder
$wsOptions.Range("B35").Value = "E:\Bot_Files\RPA WAHP Matching Automation\WABOT\"
# ZippedDirectory
$wsOptions.Range("B37").Value = "E:\Bot_Files\RPA WAHP Matching Automation\WABOT"

# Scroll the ConfigOptions sheet view up slightly (topLeftCell A30 -> A28)
$wsOptions.Activate()
$excel.ActiveWindow.ScrollRow = 28
